$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clone the row 11:13 block (values + number formats + borders/fonts/fills +
# merged-cell layout) down into the previously-empty rows 17:19, mirroring
# the "Rango" equivalence-class block used for Codigo/Descripcion/Stock, now
# reused as a template for the new "Precio" equivalence-class block.
$ws.Range("B11:G13").Copy($ws.Range("B17"))
$excel.CutCopyMode = $false

# Overwrite with the new "Precio" field content.
$ws.Range("B17").Value = "Precio"
$ws.Range("C17").Value = "Precio"
$ws.Range("D17").Value = "0<=Precio<=999999999999999"
$ws.Range("E17").Value = "CEV<07>"
$ws.Range("F17").Value = "Precio!= caracteres numéricos"
$ws.Range("G17").Value = "CENV<10>"

$ws.Range("F18").Value = "Stock>999999999999999"
$ws.Range("G18").Value = "CENV<11>"

$ws.Range("F19").Value = "Stock< 0"
$ws.Range("G19").Value = "CENV<12>"

# Match the author's final selection/cursor position.
$ws.Range("B17:G19").Select()
